$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores numeric-looking values (e.g. "304.90",
# "26.904.74") as plain text in the source data. Mark each of these
# cells as Text before writing so Excel does not silently coerce them
# into numbers and drop significant trailing/grouping digits.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Rows 12 and 13 swap: TRON and WrappedEther change places
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07511"
$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.856.99"
$ws.Range("E13").Value = "  -0.09%  "

# Price / Volume(1h) updates
$ws.Range("D2").Value = "26.904.74"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.864.53"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "304.90"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.5058"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "0.3648"
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("D9").Value = "0.07175"
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").Value = "0.8931"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("D11").Value = "20.67"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D14").Value = "94.84"
$ws.Range("E14").Value = "  +6.77%  "
$ws.Range("D15").Value = "5.231"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "0.9991"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "0.000008517"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("D18").Value = "14.22"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").Value = "0.9990"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "26.943.81"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "5.028"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "2.094.89"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D24").Value = "6.415"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "148.07"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("D26").Value = "1.780"
$ws.Range("E26").Value = "  -3.43%  "
$ws.Range("D27").Value = "17.89"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "2.080"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "113.17"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").Value = "4.706"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "4.668"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "0.09185"
$ws.Range("E32").Value = "  +1.77%  "
$ws.Range("D33").Value = "0.05134"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").Value = "0.7509"
$ws.Range("E34").Value = "  +3.65%  "
$ws.Range("D35").Value = "2.955"
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("D36").Value = "1.155"
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("D37").Value = "3.232"
$ws.Range("E37").Value = "  +6.52%  "
$ws.Range("D38").Value = "2.594"
$ws.Range("E38").Value = "  +5.61%  "
$ws.Range("D39").Value = "0.01999"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").Value = "0.5586"
$ws.Range("E40").Value = "  +6.30%  "
$ws.Range("D41").Value = "1.070"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "6.589"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "116.20"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").Value = "8.574"
$ws.Range("E44").Value = "  +4.26%  "
$ws.Range("D45").Value = "0.1474"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "0.4697"
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("D47").Value = "0.9986"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "10.06"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("D49").Value = "1.558"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").Value = "36.74"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").Value = "63.05"
$ws.Range("E51").Value = "  -1.05%  "

# Row 23 only has a Volume(1h) change
$ws.Range("E23").Value = "  -0.71%  "
